$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 435, shifting existing rows 435-523 down to 436-524.
$ws.Rows("435:435").Insert()

# Populate the newly inserted row 435 with the new data record.
$ws.Range("A435").Value = 9
$ws.Range("B435").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C435").Value = "Metropolitana"
$ws.Range("D435").Value = 44641
$ws.Range("D435").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E435").Value = 13
$ws.Range("F435").Value = 100112009
$ws.Range("G435").Value = "Acelga"
$ws.Range("H435").Value = "Sin especificar"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 43
$ws.Range("K435").Value = 18000
$ws.Range("L435").Value = 18000
$ws.Range("M435").Value = 18000
$ws.Range("N435").Value = '$/docena de atados'
$ws.Range("O435").Value = "Provincia de Cautín"
$ws.Range("P435").Value = 6000
$ws.Range("Q435").Value = 3
$ws.Range("R435").Value = "Hortaliza"
